$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on cells whose new value would otherwise be
# auto-coerced to a number by Excel, so they stay text like the rest
# of the column (matches the source inlineStr string cells).
$textCells = @('D5', 'D8', 'D10', 'D11', 'D12', 'D14', 'D15', 'D16', 'D19', 'D21', 'D22', 'D23', 'D26', 'D27', 'D28', 'D32', 'D33', 'D34', 'D35', 'D37', 'D38', 'D40', 'D41', 'D42', 'D43', 'D44', 'D46', 'D48', 'D49', 'D50', 'D51')
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply updated values scraped for this run.
$ws.Range('D2').Value = '34.784.44'
$ws.Range('E2').Value = '  -1.75%  '
$ws.Range('D3').Value = '1.876.66'
$ws.Range('E3').Value = '  -2.16%  '
$ws.Range('E4').Value = '  -1.02%  '
$ws.Range('D5').Value = '248.00'
$ws.Range('E5').Value = '  -1.96%  '
$ws.Range('E6').Value = '  -3.98%  '
$ws.Range('E7').Value = '  -0.99%  '
$ws.Range('D8').Value = '41.99'
$ws.Range('E8').Value = '  +3.04%  '
$ws.Range('E9').Value = '  -2.57%  '
$ws.Range('D10').Value = '51.29'
$ws.Range('E10').Value = '  -3.11%  '
$ws.Range('D11').Value = '0.0738'
$ws.Range('E11').Value = '  +0.44%  '
$ws.Range('D12').Value = '0.0972'
$ws.Range('E12').Value = '  -2.56%  '
$ws.Range('D13').Value = '2.149.16'
$ws.Range('E13').Value = '  -2.34%  '
$ws.Range('D14').Value = '12.87'
$ws.Range('E14').Value = '  +2.09%  '
$ws.Range('D15').Value = '0.716'
$ws.Range('E15').Value = '  -0.20%  '
$ws.Range('D16').Value = '4.91'
$ws.Range('E16').Value = '  -0.32%  '
$ws.Range('D17').Value = '1.874.22'
$ws.Range('E17').Value = '  -2.18%  '
$ws.Range('D18').Value = '34.773.90'
$ws.Range('E18').Value = '  -1.90%  '
$ws.Range('D19').Value = '72.91'
$ws.Range('E19').Value = '  -0.33%  '
$ws.Range('D20').Value = '0.0₃0821'
$ws.Range('E20').Value = '  -0.85%  '
$ws.Range('D21').Value = '245.08'
$ws.Range('D22').Value = '12.70'
$ws.Range('E22').Value = '  -3.33%  '
$ws.Range('D23').Value = '4.92'
$ws.Range('E23').Value = '  -3.05%  '
$ws.Range('E24').Value = '  -0.82%  '
$ws.Range('E25').Value = '  +3.65%  '
$ws.Range('D26').Value = '2.24'
$ws.Range('E26').Value = '  -6.67%  '
$ws.Range('D27').Value = '165.53'
$ws.Range('E27').Value = '  -1.24%  '
$ws.Range('D28').Value = '8.39'
$ws.Range('E28').Value = '  -3.62%  '
$ws.Range('E29').Value = '  -2.71%  '
$ws.Range('E30').Value = '  -4.76%  '
$ws.Range('D31').Value = '4.128.57'
$ws.Range('E31').Value = '  -0.10%  '
$ws.Range('D32').Value = '1.69'
$ws.Range('E32').Value = '  +8.57%  '
$ws.Range('D33').Value = '4.26'
$ws.Range('E33').Value = '  -2.26%  '
$ws.Range('D34').Value = '0.0580'
$ws.Range('E34').Value = '  +0.10%  '
$ws.Range('D35').Value = '4.17'
$ws.Range('E35').Value = '  -2.73%  '
$ws.Range('E36').Value = '  -1.09%  '
$ws.Range('B37').Value = 'WEMIXToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D37').Value = '1.65'
$ws.Range('E37').Value = '  -17.01%  '
$ws.Range('B38').Value = 'ImmutableX'
$ws.Range('C38').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D38').Value = '0.833'
$ws.Range('E38').Value = '  -8.93%  '
$ws.Range('E39').Value = '  -4.37%  '
$ws.Range('D40').Value = '98.41'
$ws.Range('E40').Value = '  +0.31%  '
$ws.Range('D41').Value = '16.95'
$ws.Range('E41').Value = '  -3.28%  '
$ws.Range('B42').Value = 'Kaspa'
$ws.Range('C42').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D42').Value = '0.0659'
$ws.Range('E42').Value = '  +0.97%  '
$ws.Range('B43').Value = 'VeChain'
$ws.Range('C43').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D43').Value = '0.0211'
$ws.Range('E43').Value = '  +0.05%  '
$ws.Range('D44').Value = '1.08'
$ws.Range('E44').Value = '  -4.97%  '
$ws.Range('B45').Value = 'Maker'
$ws.Range('C45').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D45').Value = '1.284.63'
$ws.Range('E45').Value = '  -4.59%  '
$ws.Range('B46').Value = 'RenderToken'
$ws.Range('C46').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D46').Value = '2.35'
$ws.Range('E46').Value = '  -5.98%  '
$ws.Range('E47').Value = '  -1.23%  '
$ws.Range('B48').Value = 'Gas'
$ws.Range('C48').Value = 'https://coinranking.com/coin/hfw0nnnLtSFc7+gas-gas'
$ws.Range('D48').Value = '12.23'
$ws.Range('E48').Value = '  +2.11%  '
$ws.Range('D49').Value = '0.0779'
$ws.Range('E49').Value = '  +7.22%  '
$ws.Range('B50').Value = 'MXToken'
$ws.Range('C50').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D50').Value = '2.73'
$ws.Range('E50').Value = '  -2.13%  '
$ws.Range('D51').Value = '6.47'
$ws.Range('E51').Value = '  -2.56%  '
